$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 1028.75
$ws.Cells.Item(32, 9).Value = 790
$ws.Cells.Item(32, 10).Value = 1108.3334
$ws.Cells.Item(32, 11).Value = 790
$ws.Cells.Item(32, 12).Value = 1108.3334
$ws.Cells.Item(32, 13).Value = -464
$ws.Cells.Item(32, 14).Value = -1760.3334
$ws.Cells.Item(74, 8).Value = 3645.15
$ws.Cells.Item(74, 9).Value = 3700.1765
$ws.Cells.Item(74, 10).Value = 3333.3333
$ws.Cells.Item(74, 11).Value = 3700.1765
$ws.Cells.Item(74, 12).Value = 3333.3333
$ws.Cells.Item(74, 13).Value = -2764.1765
$ws.Cells.Item(74, 14).Value = -5205.3333
$ws.Cells.Item(77, 8).Value = 3645.15
$ws.Cells.Item(77, 9).Value = 3700.1765
$ws.Cells.Item(77, 10).Value = 3333.3333
$ws.Cells.Item(77, 11).Value = 18500.8825
$ws.Cells.Item(77, 12).Value = 16666.6665
$ws.Cells.Item(77, 13).Value = -13820.8825
$ws.Cells.Item(77, 14).Value = -26026.6665
$ws.Cells.Item(112, 8).Value = 1307.6316
$ws.Cells.Item(112, 9).Value = 800
$ws.Cells.Item(112, 10).Value = 1335.8334
$ws.Cells.Item(112, 11).Value = 2400
$ws.Cells.Item(112, 12).Value = 4007.5002
$ws.Cells.Item(112, 13).Value = -1292
$ws.Cells.Item(112, 14).Value = -6223.5002
$ws.Cells.Item(113, 8).Value = 5557217.5
$ws.Cells.Item(113, 9).Value = 9092520
$ws.Cells.Item(113, 10).Value = 1742.2858
$ws.Cells.Item(113, 11).Value = 9092520
$ws.Cells.Item(113, 12).Value = 1742.2858
$ws.Cells.Item(113, 13).Value = -9089266
$ws.Cells.Item(113, 14).Value = -8250.2858
$ws.Cells.Item(141, 8).Value = 3059.2
$ws.Cells.Item(141, 9).Value = 1505.7142
$ws.Cells.Item(141, 10).Value = 5036.364
$ws.Cells.Item(141, 11).Value = 4517.142599999999
$ws.Cells.Item(141, 12).Value = 15109.092
$ws.Cells.Item(141, 13).Value = 662.8574000000008
$ws.Cells.Item(141, 14).Value = -25469.092
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 4344.3887
$ws.Cells.Item(88, 9).Value = 2212.5
$ws.Cells.Item(88, 10).Value = 4953.5
$ws.Cells.Item(88, 11).Value = 2212.5
$ws.Cells.Item(88, 12).Value = 4953.5
$ws.Cells.Item(88, 13).Value = -1806.5
$ws.Cells.Item(88, 14).Value = -5765.5
$ws.Cells.Item(91, 8).Value = 4344.3887
$ws.Cells.Item(91, 9).Value = 2212.5
$ws.Cells.Item(91, 10).Value = 4953.5
$ws.Cells.Item(91, 11).Value = 2212.5
$ws.Cells.Item(91, 12).Value = 4953.5
$ws.Cells.Item(91, 13).Value = -808.5
$ws.Cells.Item(91, 14).Value = -7761.5
$ws.Cells.Item(97, 8).Value = 883.3333
$ws.Cells.Item(97, 9).Value = 883.3333
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 883.3333
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = -387.3333
$ws.Cells.Item(97, 14).ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1925.3196
$ws.Cells.Item(86, 9).Value = 1934.9584
$ws.Cells.Item(86, 10).Value = 1000
$ws.Cells.Item(86, 11).Value = 1934.9584
$ws.Cells.Item(86, 12).Value = 1000
$ws.Cells.Item(86, 13).Value = -811.9584
$ws.Cells.Item(86, 14).Value = -3246
$ws.Cells.Item(89, 8).Value = 1925.3196
$ws.Cells.Item(89, 9).Value = 1934.9584
$ws.Cells.Item(89, 10).Value = 1000
$ws.Cells.Item(89, 11).Value = 9674.791999999999
$ws.Cells.Item(89, 12).Value = 5000
$ws.Cells.Item(89, 13).Value = -4058.791999999999
$ws.Cells.Item(89, 14).Value = -16232
$ws.Cells.Item(94, 8).Value = 2797.9
$ws.Cells.Item(94, 9).Value = 1897.6666
$ws.Cells.Item(94, 10).Value = 10900
$ws.Cells.Item(94, 11).Value = 1897.6666
$ws.Cells.Item(94, 12).Value = 10900
$ws.Cells.Item(94, 13).Value = -1446.6666
$ws.Cells.Item(94, 14).Value = -11802
$ws.Cells.Item(134, 8).Value = 17282268
$ws.Cells.Item(134, 9).Value = 17858190
$ws.Cells.Item(134, 10).Value = 11906994
$ws.Cells.Item(134, 11).Value = 53574570
$ws.Cells.Item(134, 12).Value = 35720982
$ws.Cells.Item(134, 13).Value = -53572035
$ws.Cells.Item(134, 14).Value = -35726052
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 3294.4443
$ws.Cells.Item(62, 9).Value = 2753.8462
$ws.Cells.Item(62, 10).Value = 4700
$ws.Cells.Item(62, 11).Value = 2753.8462
$ws.Cells.Item(62, 12).Value = 4700
$ws.Cells.Item(62, 13).Value = -2129.8462
$ws.Cells.Item(62, 14).Value = -5948
$ws.Cells.Item(65, 8).Value = 3294.4443
$ws.Cells.Item(65, 9).Value = 2753.8462
$ws.Cells.Item(65, 10).Value = 4700
$ws.Cells.Item(65, 11).Value = 13769.231
$ws.Cells.Item(65, 12).Value = 23500
$ws.Cells.Item(65, 13).Value = -10649.231
$ws.Cells.Item(65, 14).Value = -29740
$ws.Cells.Item(92, 8).Value = 19800
$ws.Cells.Item(92, 10).Value = 19800
$ws.Cells.Item(92, 12).Value = 19800
$ws.Cells.Item(92, 14).Value = -24792
$ws.Cells.Item(93, 8).Value = 12703.5
$ws.Cells.Item(93, 9).Value = 10407
$ws.Cells.Item(93, 10).Value = 15000
$ws.Cells.Item(93, 11).Value = 10407
$ws.Cells.Item(93, 12).Value = 15000
$ws.Cells.Item(93, 13).Value = -8535
$ws.Cells.Item(93, 14).Value = -18744
$ws.Cells.Item(95, 8).Value = 30000
$ws.Cells.Item(95, 10).Value = 30000
$ws.Cells.Item(95, 12).Value = 30000
$ws.Cells.Item(95, 14).Value = -35492
$ws.Cells.Item(105, 8).Value = 7966
$ws.Cells.Item(105, 9).Value = 2094.2856
$ws.Cells.Item(105, 10).Value = 21666.666
$ws.Cells.Item(105, 11).Value = 2094.2856
$ws.Cells.Item(105, 12).Value = 21666.666
$ws.Cells.Item(105, 13).Value = -347.2856000000002
$ws.Cells.Item(105, 14).Value = -25160.666
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 149.95
$ws.Cells.Item(8, 9).Value = 149.95
$ws.Cells.Item(8, 11).Value = 449.85
$ws.Cells.Item(8, 13).Value = -310.85
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 2405810
$ws.Cells.Item(70, 9).Value = 1296802.5
$ws.Cells.Item(70, 10).Value = 5055106
$ws.Cells.Item(70, 11).Value = 1296802.5
$ws.Cells.Item(70, 12).Value = 5055106
$ws.Cells.Item(70, 13).Value = -1296532.5
$ws.Cells.Item(70, 14).Value = -5055646
$ws.Cells.Item(73, 8).Value = 2405810
$ws.Cells.Item(73, 9).Value = 1296802.5
$ws.Cells.Item(73, 10).Value = 5055106
$ws.Cells.Item(73, 11).Value = 1296802.5
$ws.Cells.Item(73, 12).Value = 5055106
$ws.Cells.Item(73, 13).Value = -1295866.5
$ws.Cells.Item(73, 14).Value = -5056978
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 1735.6
$ws.Cells.Item(40, 9).Value = 1733.6471
$ws.Cells.Item(40, 10).Value = 1746.6666
$ws.Cells.Item(40, 11).Value = 1733.6471
$ws.Cells.Item(40, 12).Value = 1746.6666
$ws.Cells.Item(40, 13).Value = -1597.6471
$ws.Cells.Item(40, 14).Value = -2018.6666
$ws.Cells.Item(93, 8).Value = 30466.889
$ws.Cells.Item(93, 9).Value = 8125.5
$ws.Cells.Item(93, 10).Value = 48340
$ws.Cells.Item(93, 11).Value = 8125.5
$ws.Cells.Item(93, 12).Value = 48340
$ws.Cells.Item(93, 13).Value = -6877.5
$ws.Cells.Item(93, 14).Value = -50836
$ws.Cells.Item(136, 8).Value = 6175367
$ws.Cells.Item(136, 9).Value = 11113260
$ws.Cells.Item(136, 10).Value = 3000.625
$ws.Cells.Item(136, 11).Value = 33339780
$ws.Cells.Item(136, 12).Value = 9001.875
$ws.Cells.Item(136, 13).Value = -33337230
$ws.Cells.Item(136, 14).Value = -14101.875
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 4492.66
$ws.Cells.Item(136, 9).Value = 2709.3257
$ws.Cells.Item(136, 10).Value = 12161
$ws.Cells.Item(136, 11).Value = 8127.9771
$ws.Cells.Item(136, 12).Value = 36483
$ws.Cells.Item(136, 13).Value = -5577.9771
$ws.Cells.Item(136, 14).Value = -41583
